$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently ends at row 119 (A1:J119). Two more days of data need to
# be appended as rows 120 and 121.
#
# Copy the formatting (number formats / styles) of the last existing data row
# (119) down into the two new rows so the new cells inherit the same styles
# (date format on column B, thousands-separator number format on C:J) instead
# of picking up Excel's defaults.
$ws.Range("A119:J119").Copy($ws.Range("A120:J120"))
$ws.Range("A119:J119").Copy($ws.Range("A121:J121"))

# ---- Row 120: _id 118, 2021-04-26 ----
$ws.Range("A120").Value = 118
$ws.Range("B120").Value = 44312
$ws.Range("C120").Value = 69308
$ws.Range("D120").Value = 4696211
$ws.Range("E120").Formula = "=(D120-F120)"
$ws.Range("F120").Value = 722332
$ws.Range("G120").Value = 361166
$ws.Range("H120").Formula = "=AVERAGE(C114:C120)"
$ws.Range("I120").Formula = "=(D120-G120)"
$ws.Range("J120").Formula = "=(I120-I119)"

# ---- Row 121: _id 119, 2021-04-27 ----
$ws.Range("A121").Value = 119
$ws.Range("B121").Value = 44313
$ws.Range("C121").Value = 94819
$ws.Range("D121").Value = 4791030
$ws.Range("E121").Formula = "=(D121-F121)"
$ws.Range("F121").Value = 725126
$ws.Range("G121").Value = 362563
$ws.Range("H121").Formula = "=AVERAGE(C115:C121)"
$ws.Range("I121").Formula = "=(D121-G121)"
$ws.Range("J121").Formula = "=(I121-I120)"

# Match the workbook's new selection (the last filled-in cell, H121).
[void]$ws.Range("H121").Select()
